$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (previously "Hanoi Jones") -> "Player 2", Hit Points 50 -> 25
$ws.Range("A2").Value = "Player 2"
$ws.Range("G2").Value = 25

# Update row 3 (previously "Saoirse") -> "Player 1", Hit Points 75 -> 25
$ws.Range("A3").Value = "Player 1"
$ws.Range("G3").Value = 25

# Remove the remaining player rows (4-7): HXR-45, Caszire, Gwen, Keely
$ws.Rows("4:7").Delete()
